$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'289.78"
$ws.Range("E2").Value = "'-4.18%"
$ws.Range("D3").Value = "'30.76"
$ws.Range("E3").Value = "'-4.45%"
$ws.Range("D4").Value = "'4.889"
$ws.Range("E4").Value = "'-1.82%"
$ws.Range("D5").Value = "'0.07182"
$ws.Range("E5").Value = "'-9.22%"
$ws.Range("D6").Value = "'1.791"
$ws.Range("E6").Value = "'-16.13%"
$ws.Range("D7").Value = "'7.674"
$ws.Range("E7").Value = "'-2.05%"
$ws.Range("D8").Value = "'3.736"
$ws.Range("E8").Value = "'-1.98%"
$ws.Range("D9").Value = "'0.8993"
$ws.Range("E9").Value = "'-2.61%"
$ws.Range("D10").Value = "'0.1657"
$ws.Range("E10").Value = "'-5.22%"
$ws.Range("D11").Value = "'0.07487"
$ws.Range("E11").Value = "'-5.53%"
$ws.Range("D12").Value = "'0.08050"
$ws.Range("E12").Value = "'-7.83%"
$ws.Range("D13").Value = "'0.02988"
$ws.Range("E13").Value = "'-4.49%"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("D15").Value = "'0.001500"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("D16").Value = "'0.005697"
$ws.Range("E16").Value = "'-1.86%"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D18").Value = "'3.463"
$ws.Range("E18").Value = "'-0.14%"
$ws.Range("D19").Value = "'2.108"
$ws.Range("E19").Value = "'-7.38%"
$ws.Range("E20").Value = "'-0.29%"
$ws.Range("E21").Value = "'0.48%"
$ws.Range("D22").Value = "'4.382"
$ws.Range("E22").Value = "'1.23%"
$ws.Range("E23").Value = "'11.82%"
$ws.Range("E24").Value = "'-2.81%"
$ws.Range("D25").Value = "'0.001212"
$ws.Range("E25").Value = "'-2.11%"
$ws.Range("D26").Value = "'0.004024"
$ws.Range("E26").Value = "'-9.98%"
$ws.Range("D27").Value = "'0.0001252"
$ws.Range("E27").Value = "'0.18%"
$ws.Range("D39").Value = "'0.01644"
$ws.Range("E39").Value = "'-4.55%"
$ws.Range("D40").Value = "'0.04345"
$ws.Range("E40").Value = "'-9.46%"
$ws.Range("D41").Value = "'0.007413"
$ws.Range("E41").Value = "'-0.04%"
$ws.Range("D42").Value = "'0.1310"
$ws.Range("E42").Value = "'-3.58%"
$ws.Range("D43").Value = "'0.002043"
$ws.Range("E43").Value = "'-14.13%"
$ws.Range("D44").Value = "'0.01023"
$ws.Range("E44").Value = "'-0.20%"
$ws.Range("D45").Value = "'0.00005818"
$ws.Range("E45").Value = "'-2.55%"
$ws.Range("E46").Value = "'0.18%"
$ws.Range("D47").Value = "'2.196"
$ws.Range("E47").Value = "'167.65%"
$ws.Range("E48").Value = "'-11.41%"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "'0.18%"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("E50").Value = "'0.18%"
